# NFHS Family Aggregation Variable List — add a "Marital status 3" /
# "marital3" row, carrying the ha60 / hb60 (HIV recode) column values that
# used to be parked on the "Marital status 1" row down to their own new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("7a variables")

# Grab the HIV-recode variable names (ha60 / hb60) that are currently
# sitting in E28:F28 on the "Marital status 1" row — they belong to the
# new "Marital status 3" row instead.
$haVar = $ws.Range("E28").Value2
$hbVar = $ws.Range("F28").Value2

# Insert a fresh row right after "Marital status 2" (row 29), pushing
# "Wealth index" and everything below it down by one.
$ws.Rows.Item(30).Insert()

# Fill in the new "Marital status 3" row. (Set B before A so the shared
# string table picks up "marital3" before "Marital status 3", matching
# how the values were entered.)
$ws.Range("B30").Value = "marital3"
$ws.Range("A30").Value = "Marital status 3"
$ws.Range("E30").Value = $haVar
$ws.Range("F30").Value = $hbVar

# Clear out the old E28:F28 cells now that their values live on row 30.
$ws.Range("E28:F28").ClearContents()

# The conditional-formatting "duplicate values" rules are anchored to
# fixed row ranges, so they need to be nudged down to keep covering the
# same logical rows now that row 30 is new.
$cf = $ws.Cells.FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("B41:B48"))
$cf.Item(2).ModifyAppliesToRange($ws.Range("B49:B57"))
$cf.Item(3).ModifyAppliesToRange($ws.Range("B58:B59"))
$cf.Item(4).ModifyAppliesToRange($ws.Range("B60"))
$cf.Item(5).ModifyAppliesToRange($ws.Range("B29:B30"))
$cf.Item(6).ModifyAppliesToRange($ws.Range("B79:B88"))
$cf.Item(7).ModifyAppliesToRange($ws.Range("B78"))
$cf.Item(8).ModifyAppliesToRange($ws.Range("B90:B95"))

# Leave the view where the author ended up after this edit.
$ws.Range("C22").Select()
